$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - STL::sort
$ws.Range("B5").Value = 40000000
$ws.Range("D5").Value = 8.1
$ws.Range("F5").Value = 16.21

# Row 6 - Quick Sort
$ws.Range("B6").Value = 1000000
$ws.Range("D6").Value = 15.67
$ws.Range("F6").Value = 62.45

# Row 7 - Bubble Sort
$ws.Range("B7").Value = 10100
$ws.Range("D7").Value = 15.96
$ws.Range("F7").Value = 63.6309

# Row 8 - Insertion Sort
$ws.Range("B8").Value = 12000
$ws.Range("D8").Value = 16.31
$ws.Range("F8").Value = 65.42

# Row 9 - Selection Sort
$ws.Range("B9").Value = 13000
$ws.Range("D9").Value = 16.71
$ws.Range("F9").Value = 66.55

# Row 10 - Merge Sort
$ws.Range("B10").Value = 250000
$ws.Range("D10").Value = 13.12
$ws.Range("F10").Value = 47.36

# Update selection to F11
$ws.Range("F11").Select()
